$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Work-progress log: four new date/description rows (10-13) appended below the
# existing entries. Column A holds dd.mm.yyyy date strings that must stay plain
# text (not get auto-converted to Excel date serials), so NumberFormat is
# forced to "@" before the value is written and reset back to Normal right
# after so the cells keep the workbook's default (unstyled) look.
#
# Values are written in the same order the author's entries appear in the
# sheet's shared-string table (row 13 down to row 10), so the generated
# sharedStrings.xml ordering matches.

$ws.Range("A10:A13").NumberFormat = "@"

$ws.Range("A13").Value = "13.06.2024"
$ws.Range("B13").Value = "Koodile on tehtud viimane lihv.                                                   Viimane tegelane lisatud."

$ws.Range("A12").Value = "12.06.2024"
$ws.Range("B12").Value = "Kogu dialoogi muutmine."

$ws.Range("A11").Value = "05.06.2024"

$ws.Range("A10").Value = "03.06.2024"
$ws.Range("B10").Value = "Viimane dialoog on kirjutatud."

$ws.Range("B11").Value = "Klasside ja piltide lisamine tegelastele.                                     Mängule lõppude lisamine."

$ws.Range("A10:A13").Style = "Normal"

# Column B entries for rows 11 & 13 are long enough to wrap (same look as the
# other multi-line notes in the sheet, e.g. B9) - reuses the existing
# wrap-text style already used elsewhere in column B.
$ws.Range("B11").WrapText = $true
$ws.Range("B13").WrapText = $true

# Those two wrapped rows render at a taller row height, same as row 6/11/13 etc.
$ws.Rows.Item(11).RowHeight = 45
$ws.Rows.Item(13).RowHeight = 45

# Column B widened slightly to fit the new content.
$ws.Columns.Item(2).ColumnWidth = 26.5

# Selection follows the newly added last row, like the author leaving off at F13.
$ws.Range("F13").Select()
